$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Disinformation" + "COVID" become one paragraph: "Disinformation on COVID"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Disinformation^pCOVID", $false, $false, $false, $false, $false, $true, 1, $false, "Disinformation on COVID", 2) | Out-Null

# 2) Append new sentence fragment to the COVID detail paragraph
$r = $d.Content
$r.Find.Execute("death of many Americans.") | Out-Null
$r.InsertAfter(" According the National Institute of Health,")
# 3) New blank paragraph right after it
$r.InsertParagraphAfter()

# 4) "Climate Crisis" paragraph gets a new leading label run
$r = $d.Content
$r.Find.Execute("Climate Crisis") | Out-Null
$r.InsertBefore("Disinformation on ")

# 5) New blank paragraph right before "Self-driving cars"
$r = $d.Content
$r.Find.Execute("Self-driving cars") | Out-Null
$prevPara = $r.Paragraphs(1).Previous()
$prevPara.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# "Addictive social media algos" -> "...algorithms and mental health"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Addictive social media algos", $false, $false, $false, $false, $false, $true, 1, $false, "Addictive social media algorithms and mental health", 2) | Out-Null

# ------------------------------------------------------------------
# New sentence in the blank paragraph right after "Medicine"
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Medicine") | Out-Null
$medNext = $r.Paragraphs(1).Next()
$medNext.Range.InsertBefore("With increased automation in the drug development process, human oversite is required to ensure")
$medNext.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# New sentence in the blank paragraph right after "AI-weaponry"
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("AI-weaponry") | Out-Null
$aiNext = $r.Paragraphs(1).Next()
$aiNext.Range.InsertBefore("Besides malicious human actors, the release of a sufficiently intelligent self-preserving program onto the internet can pose a threat to internet infrastructure and communications networks.")
$aiNext.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# "Job loss" -> "Job loss from automation" + new blank paragraph after
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Job loss") | Out-Null
$r.InsertAfter(" from automation")
$r.InsertParagraphAfter()

# ------------------------------------------------------------------
# "Explicitly define privacy, ... responsibility in developers. " rewritten
# ------------------------------------------------------------------
$oldExplicit = "Explicitly define privacy, safety, security, transparency, responsibility in developers. "
$newExplicit = "Explicitly define privacy, safety, security, transparency, and responsibility for developers to be held accountable to. "
$d.Content.Find.Execute($oldExplicit, $false, $false, $false, $false, $false, $true, 1, $false, $newExplicit, 2) | Out-Null

# ------------------------------------------------------------------
# Drop the "Twitter bot becoming a Nazi" aside, leaving just a space
# ------------------------------------------------------------------
$d.Content.Find.Execute("Twitter bot becoming a Nazi", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ------------------------------------------------------------------
# "Between both of the previous presidential administrations, the NSTC have been tasked"
#   -> "Between both previous presidential administrations, the NSTC has been tasked"
# ------------------------------------------------------------------
$oldBetween = "Between both of the previous presidential administrations, the NSTC have been tasked"
$newBetween = "Between both previous presidential administrations, the NSTC has been tasked"
$d.Content.Find.Execute($oldBetween, $false, $false, $false, $false, $false, $true, 1, $false, $newBetween, 2) | Out-Null

# ------------------------------------------------------------------
# "National Science and Technology Council (NSTC)" gains a trailing clause
# ------------------------------------------------------------------
$oldNstc = "National Science and Technology Council (NSTC)"
$newNstc = "National Science and Technology Council" + [char]0x2019 + "s (NSTC) Special Committee on Artificial Intelligence to establish norms and practices."
$d.Content.Find.Execute($oldNstc, $false, $false, $false, $false, $false, $true, 1, $false, $newNstc, 2) | Out-Null

# ------------------------------------------------------------------
# Consulting section gets rewritten
# ------------------------------------------------------------------
$oldConsult = "On regular, recurring basis, the Agency will consult with experts in industry, government, and academia in order to keep up to date on the latest changes to the practice."
$newConsult = "On a regular, recurring basis, directors of the Agency will consult with experts in industry, government, and academia in order to keep up to date on the latest changes to the development of Artificial Intelligence systems."
$d.Content.Find.Execute($oldConsult, $false, $false, $false, $false, $false, $true, 1, $false, $newConsult, 2) | Out-Null

Write-Output "All edits applied"
